$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("D:D").Insert()
$ws.Range("E7").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("E38").Copy()
$ws.Range("D38").PasteSpecial(-4122)
$ws.Range("E80").Copy()
$ws.Range("D80").PasteSpecial(-4122)
$ws.Range("E8").Copy()
$ws.Range("D8:D102").PasteSpecial(-4122)
